# Auto-generated Excel COM-interop script
# Applies numeric corrections to columns H-N across several rows/sheets
# per the commit 'chore: update Sheets via scheduled runner'.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1037.3
$ws.Range("I19").Value = 499.5
$ws.Range("K19").Value = 499.5
$ws.Range("M19").Value = -324.5

$ws.Range("H40").Value = 10962.889
$ws.Range("J40").Value = 11494.833
$ws.Range("L40").Value = 11494.833
$ws.Range("N40").Value = -11844.833

$ws.Range("H100").Value = 3935
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws.Range("H103").Value = 896.8570999999999
$ws.Range("J103").Value = 450
$ws.Range("L103").Value = 1350
$ws.Range("N103").Value = -2522

$ws.Range("H137").Value = 2916.6667

$ws.Range("H138").Value = 2318.9092
$ws.Range("J138").Value = 2433.16
$ws.Range("L138").Value = 7299.48
$ws.Range("N138").Value = -17579.48


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H51").Value = 54000
$ws.Range("J51").Value = 54000
$ws.Range("L51").Value = 54000
$ws.Range("N51").Value = -55512

$ws.Range("H61").Value = 4462.778
$ws.Range("I61").Value = 4038.0344
$ws.Range("J61").Value = 6222.4287
$ws.Range("K61").Value = 4038.0344
$ws.Range("L61").Value = 6222.4287
$ws.Range("M61").Value = -3826.0344
$ws.Range("N61").Value = -6646.4287

$ws.Range("H74").Value = 7869.5
$ws.Range("I74").Value = 5749
$ws.Range("J74").Value = 8293.6
$ws.Range("K74").Value = 5749
$ws.Range("L74").Value = 8293.6
$ws.Range("M74").Value = -4875
$ws.Range("N74").Value = -10041.6

$ws.Range("H77").Value = 7869.5
$ws.Range("I77").Value = 5749
$ws.Range("J77").Value = 8293.6
$ws.Range("K77").Value = 28745
$ws.Range("L77").Value = 41468
$ws.Range("M77").Value = -24377
$ws.Range("N77").Value = -50204

$ws.Range("H122").Value = 2301.7878
$ws.Range("I122").Value = 1945.1072
$ws.Range("K122").Value = 5835.321599999999
$ws.Range("M122").Value = -3385.321599999999

$ws.Range("H136").Value = 4462.778
$ws.Range("I136").Value = 4038.0344
$ws.Range("J136").Value = 6222.4287
$ws.Range("K136").Value = 12114.1032
$ws.Range("L136").Value = 18667.2861
$ws.Range("M136").Value = -9564.1032
$ws.Range("N136").Value = -23767.2861


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 16699.75
$ws.Range("J74").Value = 12000
$ws.Range("L74").Value = 12000
$ws.Range("N74").Value = -13872

$ws.Range("H77").Value = 16699.75
$ws.Range("J77").Value = 12000
$ws.Range("L77").Value = 36000
$ws.Range("N77").Value = -45360

$ws.Range("H94").Value = 9979.5
$ws.Range("I94").Value = 9979.5
$ws.Range("K94").Value = 9979.5
$ws.Range("M94").Value = -9528.5

$ws.Range("H120").Value = 38999
$ws.Range("J120").Value = 38999
$ws.Range("L120").Value = 38999
$ws.Range("N120").Value = -48675

$ws.Range("H141").Value = 53259.668
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 53259.668
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 53259.668
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -63619.668


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4777.5
$ws.Range("I86").Value = 4630.3335
$ws.Range("K86").Value = 4630.3335
$ws.Range("M86").Value = -3507.3335

$ws.Range("H89").Value = 4777.5
$ws.Range("I89").Value = 4630.3335
$ws.Range("K89").Value = 23151.6675
$ws.Range("M89").Value = -17535.6675

$ws.Range("H107").Value = 1598.8667
$ws.Range("I107").Value = 1284.5
$ws.Range("K107").Value = 1284.5
$ws.Range("M107").Value = 635.5

$ws.Range("H134").Value = 3190.4473
$ws.Range("I134").Value = 2889.0908
$ws.Range("K134").Value = 8667.2724
$ws.Range("M134").Value = -6132.2724

$ws.Range("H140").Value = 91926.336
$ws.Range("J140").Value = 112889.5
$ws.Range("L140").Value = 112889.5
$ws.Range("N140").Value = -123249.5


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1011.2143
$ws.Range("J5").Value = 1342.5714
$ws.Range("L5").Value = 4027.7142
$ws.Range("N5").Value = -4251.7142

$ws.Range("H55").Value = 571.4286
$ws.Range("I55").Value = 571.4286
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 1714.2858
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -1537.2858
$ws.Range("N55").ClearContents()

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H135").Value = 1011.2143
$ws.Range("J135").Value = 1342.5714
$ws.Range("L135").Value = 12083.1426
$ws.Range("N135").Value = -17153.1426


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 848
$ws.Range("I2").Value = 921.45
$ws.Range("K2").Value = 921.45
$ws.Range("M2").Value = -808.45

$ws.Range("H15").Value = 59999
$ws.Range("J15").Value = 59999
$ws.Range("L15").Value = 59999
$ws.Range("N15").Value = -60575

$ws.Range("H26").Value = 33332.668
$ws.Range("J26").Value = 39999.5
$ws.Range("L26").Value = 39999.5
$ws.Range("N26").Value = -40559.5

$ws.Range("H46").Value = 29497.875
$ws.Range("J46").Value = 29497.875
$ws.Range("L46").Value = 29497.875
$ws.Range("N46").Value = -29809.875

$ws.Range("H50").Value = 33332.668
$ws.Range("J50").Value = 39999.5
$ws.Range("L50").Value = 39999.5
$ws.Range("N50").Value = -40995.5

$ws.Range("H81").Value = 59999
$ws.Range("J81").Value = 59999
$ws.Range("L81").Value = 59999
$ws.Range("N81").Value = -61995

$ws.Range("H84").Value = 59999
$ws.Range("J84").Value = 59999
$ws.Range("L84").Value = 179997
$ws.Range("N84").Value = -189981

$ws.Range("H109").Value = 40000
$ws.Range("J109").Value = 40000
$ws.Range("L109").Value = 40000
$ws.Range("N109").Value = -42080

$ws.Range("H132").Value = 4949.317
$ws.Range("I132").Value = 3742.476
$ws.Range("J132").Value = 6216.5
$ws.Range("K132").Value = 11227.428
$ws.Range("L132").Value = 18649.5
$ws.Range("M132").Value = -8697.428
$ws.Range("N132").Value = -23709.5


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3532.2307
$ws.Range("J46").Value = 3897.7
$ws.Range("L46").Value = 3897.7
$ws.Range("N46").Value = -4273.7

$ws.Range("H132").Value = 4220.696
$ws.Range("I132").Value = 2228.8
$ws.Range("J132").Value = 17500
$ws.Range("K132").Value = 6686.400000000001
$ws.Range("L132").Value = 52500
$ws.Range("M132").Value = -4156.400000000001
$ws.Range("N132").Value = -57560

